$p = $ppt.ActivePresentation
$m1 = $p.SlideMaster.Theme.ThemeColorScheme
$m2 = $p.NotesMaster.Theme.ThemeColorScheme
for ($i=1; $i -le 12; $i++) {
    $c1 = $m1.Colors($i).RGB
    $c2 = $m2.Colors($i).RGB
    Write-Output "$i : master=$c1 notes=$c2 same=$($c1 -eq $c2)"
}
